# This script applies the same set of cell edits to both the "展览" and
# "全部类型" worksheets, which (in the source workbook) contain identical
# event listing data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # "Sold out" -> "Not available for sale"
    $ws.Range("G3").Value = "不可售"

    # "Interested" counter updates
    $ws.Range("F5").Value = 125
    $ws.Range("F10").Value = 1171
    $ws.Range("F15").Value = 122
    $ws.Range("F21").Value = 307
    $ws.Range("F22").Value = 1689
    $ws.Range("F27").Value = 303
    $ws.Range("F28").Value = 208
    $ws.Range("F29").Value = 4033
    $ws.Range("F31").Value = 473
    $ws.Range("F33").Value = 1040
    $ws.Range("F36").Value = 273
    $ws.Range("F38").Value = 154
}
